# Hortaliza, Vega Modelo de Temuco - Perejil
# A new weekly observation is inserted at row 326; every existing data row
# from 326..352 shifts down by one (327..353), carrying along the columns
# that vary per-observation (Fecha, Volumen, Precio min/max/prom, Origen,
# Precio $/Kg). The static/category columns (A,B,C,E,F,G,H,I,N,Q,R) are
# identical across the whole block, so only the moving columns need to be
# rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 326
$lastOldRow = 352
$newLastRow = 353

# Columns that carry per-row observation data and shift down by one row.
$movingCols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the moving columns across the whole
# block (326..352) before writing anything, so the shift-down can be
# applied safely regardless of write order.
$snapshot = @{}
foreach ($col in $movingCols) {
    $colVals = @{}
    for ($r = $firstRow; $r -le $lastOldRow; $r++) {
        $colVals[$r] = $ws.Range("$col$r").Value2
    }
    $snapshot[$col] = $colVals
}

# New row 353 is a full copy of the old row 352 (static columns already
# match; copy them explicitly too so the appended row is complete).
$staticCols = @("A", "B", "C", "E", "F", "G", "H", "I", "N", "Q", "R")
foreach ($col in $staticCols) {
    $ws.Range("$col$newLastRow").Value2 = $ws.Range("$col$lastOldRow").Value2
}

# Shift the moving-column values down by one row: new row r (327..353)
# gets the old value that was in row r-1.
for ($r = $newLastRow; $r -ge ($firstRow + 1); $r--) {
    foreach ($col in $movingCols) {
        $ws.Range("$col$r").Value2 = $snapshot[$col][$r - 1]
    }
    $ws.Range("D$r").NumberFormat = $ws.Range("D$firstRow").NumberFormat
}

# Row 326 becomes the new observation.
$ws.Range("D$firstRow").Value2 = 44783
$ws.Range("J$firstRow").Value2 = 40
$ws.Range("K$firstRow").Value2 = 3600
$ws.Range("L$firstRow").Value2 = 3600
$ws.Range("M$firstRow").Value2 = 3600
$ws.Range("O$firstRow").Value2 = "Región Metropolitana"
$ws.Range("P$firstRow").Value2 = 1200

$wb.Save()
